$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values for the specified rows to reflect repulled data
$ws.Range("F3").Value = -6
$ws.Range("F6").Value = -7
$ws.Range("F7").Value = -3
$ws.Range("F13").Value = -1
$ws.Range("F15").Value = 8
$ws.Range("F18").Value = -1
